$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25)

$vals = @(
    @(1.02,1.018623176781927,1.024047941851572,1.019854250497469,1.032164508047733,1.028596566759497,1.023830561357466,1.026877588401623,1.022696264795027,1.034970528133968,1.011980991928231),
    @(1.02,1.019700392141955,1.024826104283549,1.02077192242974,1.034336153691982,1.028828947857333,1.024543509490078,1.027462837790671,1.023419744398527,1.036947245532658,1.012221928229621),
    @(1.02,1.020397300185583,1.025329252403386,1.021365971661623,1.035735540219732,1.028977463979154,1.02500419355662,1.027840505261647,1.023887516494799,1.038220173893348,1.012377479995612),
    @(1.02,1.020690252513746,1.025540686796589,1.02161577136763,1.036322482954467,1.029039458876218,1.025197712547429,1.02799903148241,1.024084080563763,1.038753873196995,1.012442790333952),
    @(1.02,1.020739438910161,1.025576182357995,1.021657717454723,1.03642095444111,1.029049842261565,1.025230196314668,1.028025634381187,1.02411707943481,1.038843400096976,1.012453751329014),
    @(1.02,1.020401214736409,1.025332077949283,1.021369309255548,1.035743388279291,1.028978294089985,1.025006779965367,1.027842624459066,1.023890143335804,1.038227310830315,1.012378353003899),
    @(1.02,1.018987252976018,1.024311003566613,1.020164329753376,1.032899649188759,1.02867548506945,1.024071638784455,1.027075589594952,1.022940844640164,1.03563985862222,1.012062490098017),
    @(1.02,1.016494667706385,1.022508834850264,1.018042921338489,1.02784262633696,1.028127654655971,1.022418850294671,1.025716059888531,1.021265212096009,1.031032046176218,1.011503208326929),
    @(1.02,1.014832171625443,1.021305385357996,1.016629898089391,1.02443828580348,1.027752753385907,1.021313602922345,1.02480432094096,1.020146161959072,1.027925708076939,1.011128530103709),
    @(1.02,1.014112086380986,1.020783790714897,1.016018329254592,1.022955877590179,1.027588096590122,1.020834200514592,1.024408235325696,1.019661123961071,1.026572021169985,1.010965854109432),
    @(1.02,1.013844580499849,1.020589971963041,1.015791206346168,1.022403956779492,1.02752658469152,1.020656004117533,1.02426091527386,1.019480885744787,1.026067868503834,1.010905362863421),
    @(1.02,1.01390196299905,1.020631550169345,1.015839923142829,1.022522404626768,1.027539795122845,1.02069423352799,1.024292524827432,1.01951955079843,1.026176071899135,1.010918341439774),
    @(1.02,1.014089974970428,1.020767771139111,1.015999554373969,1.022910282163377,1.027583019171069,1.020819473302648,1.024396061815355,1.019646226919915,1.02653037516828,1.010960855235472),
    @(1.02,1.014205810761742,1.020851691362957,1.016097913844352,1.023149094367293,1.027609604352818,1.020896621056979,1.024459828327017,1.019724266485184,1.026748495320531,1.010987040613292),
    @(1.02,1.014879956636432,1.021339991386672,1.016670491702875,1.024536489912218,1.027763632007691,1.021345401826931,1.02483058040025,1.020178342086515,1.028015362958671,1.011139317118479),
    @(1.019999999999999,1.015302771915246,1.021646156246189,1.017029728804213,1.025404514305107,1.027859626284448,1.021626688861313,1.025062795246819,1.020463041921634,1.028807700783044,1.01123471867335),
    @(1.02,1.015549372460862,1.021824689538618,1.017239292771852,1.025910019579038,1.02791539422479,1.021790679515402,1.025198117229768,1.02062905609907,1.029269028712133,1.011290322564939),
    @(1.02,1.015633453467254,1.021885556757764,1.01731075328616,1.026082249560185,1.027934371739236,1.02184658267757,1.025244237367289,1.020685654802528,1.029426190078919,1.01130927490406),
    @(1.02,1.015257409968144,1.021613312575174,1.016991183273179,1.025311466403529,1.027849350181403,1.021596517638292,1.025037893731922,1.020432501133211,1.028722776484462,1.011224487359144),
    @(1.02,1.014034611090816,1.020727659546568,1.015952545831936,1.022796097790999,1.027570300474064,1.02078259677206,1.02436557818112,1.019608925997642,1.026426078791905,1.010948337810929),
    @(1.02,1.013265590211567,1.020170378183225,1.015299749864645,1.021207110880638,1.027392818771223,1.020270128532292,1.023941730724577,1.019090685763015,1.024974321860671,1.010774328788072),
    @(1.02,1.013673282267789,1.02046584532969,1.015645787274745,1.022050185131939,1.027487098532466,1.020541866716482,1.024166528439564,1.01936545545925,1.025744671343255,1.010866610673895),
    @(1.02,1.015277907143123,1.021628153363481,1.017008600256356,1.025353513211312,1.027853994201918,1.021610150960681,1.025049146041173,1.020446301344669,1.028761152680546,1.01122911058022),
    @(1.02,1.017139189676581,1.022975088012227,1.018591132984799,1.029155649167152,1.028270980436014,1.022846728464215,1.026068475440108,1.021698745473913,1.032229207218026,1.011648116419564)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $addr = "$($cols[$j])$($rows[$i])"
        $ws.Range($addr).Value = $vals[$i][$j]
    }
}

Write-Output "applied $($rows.Length) rows"